$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161185741424561
$ws.Range("B1").Value = 2.413383960723877
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.380521297454834
$ws.Range("E1").Value = 1.230362176895142
